# Apply "Optuna Attempt (go back with original)" changes to the
# "Forecast Comparison" worksheet (sheet1.xml).
#
# Updates the Seasonality Index (column L) values for rows 2-17,
# and for row 5 also updates Inventory Coverage (H5) and Reorder
# Urgency (J5); for row 6 also updates Inventory Coverage (H6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Inventory Coverage (column H)
$ws.Range("H5").Value = 1.1
$ws.Range("H6").Value = 0.1

# Reorder Urgency (column J)
$ws.Range("J5").Value = "Normal"

# Seasonality Index (column L)
$ws.Range("L2").Value = 0.9399999999999999
$ws.Range("L3").Value = 0.88
$ws.Range("L4").Value = 1.11
$ws.Range("L5").Value = 0.9399999999999999
$ws.Range("L6").Value = 0.9399999999999999
$ws.Range("L7").Value = 1.08
$ws.Range("L8").Value = 0.82
$ws.Range("L9").Value = 1.15
$ws.Range("L10").Value = 1
$ws.Range("L11").Value = 1.17
$ws.Range("L12").Value = 1.18
$ws.Range("L13").Value = 0.92
$ws.Range("L14").Value = 0.91
$ws.Range("L15").Value = 0.88
$ws.Range("L16").Value = 0.9399999999999999
$ws.Range("L17").Value = 0.91
